# Updated cryptos list values to match the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.294.09'
$ws.Range('E2').Value = '  +0.97%  '
$ws.Range('D3').Value = '1.564.85'
$ws.Range('E3').Value = '  +0.13%  '
$ws.Range('E4').Value = '  -0.34%  '
$ws.Range('D5').Value = '''211.32'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +1.48%  '
$ws.Range('E6').Value = '  -0.22%  '
$ws.Range('E7').Value = '  -0.46%  '
$ws.Range('D8').Value = '''22.24'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  +0.86%  '
$ws.Range('E9').Value = '  +0.29%  '
$ws.Range('E10').Value = '  -0.48%  '
$ws.Range('E11').Value = '  +2.21%  '
$ws.Range('D12').Value = '1.787.38'
$ws.Range('E12').Value = '  +0.13%  '
$ws.Range('D13').Value = '1.572.87'
$ws.Range('E13').Value = '  +0.43%  '
$ws.Range('E15').Value = '  -0.29%  '
$ws.Range('D16').Value = '27.293.58'
$ws.Range('E16').Value = '  +0.94%  '
$ws.Range('D17').Value = '''61.83'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  -0.20%  '
$ws.Range('D18').Value = '''218.26'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  +1.10%  '
$ws.Range('B19').Value = 'ShibaInu'
$ws.Range('C19').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D19').Value = '0.0₃0707'
$ws.Range('E19').Value = '  -0.28%  '
$ws.Range('B20').Value = 'Chainlink'
$ws.Range('C20').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D20').Value = '''7.45'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +1.09%  '
$ws.Range('E21').Value = '  -0.37%  '
$ws.Range('E22').Value = '  +0.02%  '
$ws.Range('D23').Value = '''9.37'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +1.80%  '
$ws.Range('D24').Value = '''1.95'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  -0.10%  '
$ws.Range('D25').Value = '''151.40'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -1.29%  '
$ws.Range('E26').Value = '  +0.52%  '
$ws.Range('E27').Value = '  +1.25%  '
$ws.Range('E28').Value = '  -0.17%  '
$ws.Range('E29').Value = '  -0.43%  '
$ws.Range('E30').Value = '  +2.12%  '
$ws.Range('E31').Value = '  -0.44%  '
$ws.Range('E32').Value = '  +0.41%  '
$ws.Range('D33').Value = '1.460.25'
$ws.Range('E33').Value = '  +2.17%  '
$ws.Range('E35').Value = '  +5.26%  '
$ws.Range('E36').Value = '  +1.11%  '
$ws.Range('E37').Value = '  +0.28%  '
$ws.Range('E39').Value = '  +1.12%  '
$ws.Range('E40').Value = '  -0.62%  '
$ws.Range('E41').Value = '  +0.81%  '
$ws.Range('E42').Value = '  -0.42%  '
$ws.Range('E43').Value = '  +1.66%  '
$ws.Range('D44').Value = '''0.975'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -2.63%  '
$ws.Range('D45').Value = '''64.44'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -0.41%  '
$ws.Range('D46').Value = '''1.76'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  +0.19%  '
$ws.Range('D47').Value = '1.701.96'
$ws.Range('E47').Value = '  -0.09%  '
$ws.Range('D48').Value = '''85.94'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -1.29%  '
$ws.Range('D49').Value = '0.0₆0104'
$ws.Range('E49').Value = '  -0.06%  '
$ws.Range('D50').Value = '''0.0525'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +1.31%  '
$ws.Range('E51').Value = '  -1.27%  '
